# Auto-generated edit script: updates Aegis_Profits market-data values per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 811.875
$ws.Range("I70").Value = 734
$ws.Range("J70").Value = 858.6
$ws.Range("K70").Value = 2202
$ws.Range("L70").Value = 2575.8
$ws.Range("M70").Value = -1932
$ws.Range("N70").Value = -3115.8
$ws.Range("H73").Value = 811.875
$ws.Range("I73").Value = 734
$ws.Range("J73").Value = 858.6
$ws.Range("K73").Value = 2202
$ws.Range("L73").Value = 2575.8
$ws.Range("M73").Value = -1266
$ws.Range("N73").Value = -4447.8
$ws.Range("H76").Value = 3445.75
$ws.Range("I76").Value = 2943.625
$ws.Range("K76").Value = 2943.625
$ws.Range("M76").Value = -2628.625
$ws.Range("H79").Value = 3445.75
$ws.Range("I79").Value = 2943.625
$ws.Range("K79").Value = 2943.625
$ws.Range("M79").Value = -1851.625
$ws.Range("H86").Value = 6882.6665
$ws.Range("I86").Value = 5750.0527
$ws.Range("J86").Value = 9572.625
$ws.Range("K86").Value = 5750.0527
$ws.Range("L86").Value = 9572.625
$ws.Range("M86").Value = -4627.0527
$ws.Range("N86").Value = -11818.625
$ws.Range("H89").Value = 6882.6665
$ws.Range("I89").Value = 5750.0527
$ws.Range("J89").Value = 9572.625
$ws.Range("K89").Value = 28750.2635
$ws.Range("L89").Value = 47863.125
$ws.Range("M89").Value = -23134.2635
$ws.Range("N89").Value = -59095.125
$ws.Range("H111").Value = 5578.36
$ws.Range("I111").Value = 9441.416999999999
$ws.Range("J111").Value = 2012.4615
$ws.Range("K111").Value = 28324.251
$ws.Range("L111").Value = 6037.3845
$ws.Range("M111").Value = -25257.251
$ws.Range("N111").Value = -12171.3845
$ws.Range("H129").Value = 976.1778
$ws.Range("J129").Value = 993.36584
$ws.Range("L129").Value = 2980.09752
$ws.Range("N129").Value = -12980.09752
$ws.Range("H132").Value = 9267024
$ws.Range("I132").Value = 9267024
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 27801072
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -27798542
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2250.2
$ws.Range("I45").Value = 2000.5
$ws.Range("J45").Value = 2568
$ws.Range("K45").Value = 2000.5
$ws.Range("L45").Value = 2568
$ws.Range("M45").Value = -1623.5
$ws.Range("N45").Value = -3322
$ws.Range("H74").Value = 1988.6976
$ws.Range("I74").Value = 1440.625
$ws.Range("J74").Value = 2681
$ws.Range("K74").Value = 1440.625
$ws.Range("L74").Value = 2681
$ws.Range("M74").Value = -566.625
$ws.Range("N74").Value = -4429
$ws.Range("H77").Value = 1988.6976
$ws.Range("I77").Value = 1440.625
$ws.Range("J77").Value = 2681
$ws.Range("K77").Value = 7203.125
$ws.Range("L77").Value = 13405
$ws.Range("M77").Value = -2835.125
$ws.Range("N77").Value = -22141

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2292.4285
$ws.Range("I134").Value = 2328.348
$ws.Range("K134").Value = 6985.044
$ws.Range("M134").Value = -4450.044

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 5000
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H63").Value = 43000
$ws.Range("J63").Value = 43000
$ws.Range("L63").Value = 43000
$ws.Range("N63").Value = -44372
$ws.Range("H64").Value = 47598
$ws.Range("J64").Value = 47598
$ws.Range("L64").Value = 47598
$ws.Range("N64").Value = -48094
$ws.Range("H66").Value = 43000
$ws.Range("J66").Value = 43000
$ws.Range("L66").Value = 129000
$ws.Range("N66").Value = -135864
$ws.Range("H67").Value = 47598
$ws.Range("J67").Value = 47598
$ws.Range("L67").Value = 47598
$ws.Range("N67").Value = -49314
$ws.Range("H69").Value = 16500
$ws.Range("I69").Value = 16500
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 16500
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -15751
$ws.Range("N69").ClearContents()
$ws.Range("H70").Value = 16499.666
$ws.Range("J70").Value = 16499.666
$ws.Range("L70").Value = 16499.666
$ws.Range("N70").Value = -17129.666
$ws.Range("H72").Value = 16500
$ws.Range("I72").Value = 16500
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 49500
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -45756
$ws.Range("N72").ClearContents()
$ws.Range("H73").Value = 16499.666
$ws.Range("J73").Value = 16499.666
$ws.Range("L73").Value = 16499.666
$ws.Range("N73").Value = -18683.666
$ws.Range("H134").Value = 1428.862
$ws.Range("I134").Value = 1159.4546
$ws.Range("J134").Value = 2275.5715
$ws.Range("K134").Value = 3478.3638
$ws.Range("L134").Value = 6826.7145
$ws.Range("M134").Value = -943.3638000000001
$ws.Range("N134").Value = -11896.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5518.951
$ws.Range("I5").Value = 529.7646999999999
$ws.Range("J5").Value = 29752.143
$ws.Range("K5").Value = 1589.2941
$ws.Range("L5").Value = 89256.429
$ws.Range("M5").Value = -1477.2941
$ws.Range("N5").Value = -89480.429
$ws.Range("H131").Value = 1358.7843
$ws.Range("J131").Value = 1395.6364
$ws.Range("L131").Value = 4186.9092
$ws.Range("N131").Value = -14266.9092
$ws.Range("H132").Value = 2488.611
$ws.Range("I132").Value = 1533.3334
$ws.Range("J132").Value = 2679.6667
$ws.Range("K132").Value = 13800.0006
$ws.Range("L132").Value = 24117.0003
$ws.Range("M132").Value = -11270.0006
$ws.Range("N132").Value = -29177.0003
$ws.Range("H135").Value = 5518.951
$ws.Range("I135").Value = 529.7646999999999
$ws.Range("J135").Value = 29752.143
$ws.Range("K135").Value = 4767.882299999999
$ws.Range("L135").Value = 267769.287
$ws.Range("M135").Value = -2232.882299999999
$ws.Range("N135").Value = -272839.287

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2040.7778
$ws.Range("I7").Value = 1472.6666
$ws.Range("J7").Value = 2608.889
$ws.Range("K7").Value = 1472.6666
$ws.Range("L7").Value = 2608.889
$ws.Range("M7").Value = -1360.6666
$ws.Range("N7").Value = -2832.889
$ws.Range("H126").Value = 2040.7778
$ws.Range("I126").Value = 1472.6666
$ws.Range("J126").Value = 2608.889
$ws.Range("K126").Value = 4417.9998
$ws.Range("L126").Value = 7826.667
$ws.Range("M126").Value = -1947.9998
$ws.Range("N126").Value = -12766.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 40388.332
$ws.Range("J105").Value = 40388.332
$ws.Range("L105").Value = 40388.332
$ws.Range("N105").Value = -47376.332
$ws.Range("H136").Value = 844.0645
$ws.Range("I136").Value = 446.64
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 1339.92
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = 1210.08
$ws.Range("N136").Value = -12600
